$wb = $excel.ActiveWorkbook

# --- Sheet "Schedule": update Cost ($) and Unit Cost ($/ML) for rows 3 and 4 ---
$wsSchedule = $wb.Worksheets.Item("Schedule")
$wsSchedule.Range("E3").Value = 734.7728894999999
$wsSchedule.Range("F3").Value = 27.76919461451247
$wsSchedule.Range("E4").Value = 299.78164125
$wsSchedule.Range("F4").Value = 8.81192361111111

# --- Sheet "Detailed": update Price (col B) and Type (col C) values ---
$wsDetailed = $wb.Worksheets.Item("Detailed")
$wsDetailed.Range("B35").Value = -10
$wsDetailed.Range("B36").Value = -11.01
$wsDetailed.Range("B37").Value = 3.40587
$wsDetailed.Range("C37").Value = "historical"
$wsDetailed.Range("B38").Value = 11.21034
$wsDetailed.Range("C38").Value = "historical"
$wsDetailed.Range("B39").Value = 10.51865
$wsDetailed.Range("B40").Value = 36.25
$wsDetailed.Range("B43").Value = 53.50014
$wsDetailed.Range("B44").Value = 52.71663
$wsDetailed.Range("B45").Value = 57.03541
$wsDetailed.Range("B46").Value = 42.85901
$wsDetailed.Range("B47").Value = 57.03043
$wsDetailed.Range("B50").Value = 56.98
$wsDetailed.Range("B51").Value = 56.98
$wsDetailed.Range("B54").Value = 42.03323
$wsDetailed.Range("B55").Value = 49.08837
$wsDetailed.Range("B56").Value = 50.00473
$wsDetailed.Range("B61").Value = 58.01162
$wsDetailed.Range("B62").Value = 59.80935
$wsDetailed.Range("B64").Value = 26.35766
$wsDetailed.Range("B65").Value = 22.07
$wsDetailed.Range("B69").Value = 0.7
$wsDetailed.Range("B70").Value = 0.59
$wsDetailed.Range("B71").Value = 22.07
$wsDetailed.Range("B72").Value = 22.07
$wsDetailed.Range("B73").Value = 23.53344
$wsDetailed.Range("B74").Value = 23.56806
$wsDetailed.Range("B75").Value = 26.44148
$wsDetailed.Range("B76").Value = 26.60196
$wsDetailed.Range("B77").Value = 22.93241
$wsDetailed.Range("B78").Value = 27.00337
$wsDetailed.Range("B79").Value = 29.7184
$wsDetailed.Range("B80").Value = 31.71157
$wsDetailed.Range("B81").Value = 0
$wsDetailed.Range("B82").Value = -5.00046
$wsDetailed.Range("B83").Value = -6
$wsDetailed.Range("B85").Value = -4.04154
$wsDetailed.Range("B86").Value = 12.20887
$wsDetailed.Range("B87").Value = 25.73188
$wsDetailed.Range("B90").Value = 57.09
$wsDetailed.Range("B92").Value = 55.133
$wsDetailed.Range("B93").Value = 56.24403
